$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Website" key/value pair (placed at row 3 in the final layout) ---
$ws.Range("A3").Value = "Website"
$ws.Range("B3").Value = "https://www.ilabquality.com/"

# --- Add new "Text to Validate" key/value pair (placed at row 7 in the final layout) ---
$ws.Range("A7").Value = "Text to Validate"
$ws.Range("B7").Value = "Please complete this required field."
# Style this value cell with a custom font color (RGB 32,33,36 -> #202124)
$ws.Range("B7").Font.Color = 2367776

# --- Add new "Broswer Type" key/value pair (placed at row 2 in the final layout) ---
$ws.Range("A2").Value = "Broswer Type"
$ws.Range("B2").Value = "chrome"

# --- Re-write the pre-existing rows (First Name / Last Name / Email) into their
#     new, shifted-down positions (rows 4-6) ---
$ws.Range("A4").Value = "First Name"
$ws.Range("B4").Value = "Test"
$ws.Range("A5").Value = "Last Name"
$ws.Range("B5").Value = "Automation"
$ws.Range("A6").Value = "Email"
$ws.Range("B6").Value = "automationAssessment@iLABQuality.com"

# --- Grow the Excel table so it covers the new data range ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B7"))

# --- Widen column A to fit the longer labels ---
$ws.Columns.Item(1).ColumnWidth = 23.29

# --- Update the selected cell ---
$ws.Range("B2").Select() | Out-Null

# --- Switch the page to portrait orientation ---
$ws.PageSetup.Orientation = 1
